# Update cryptos list (simulating a scheduled GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.010.63"
$ws.Range("E2").Value = "  +3.45%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.471.34"
$ws.Range("E3").Value = "  +3.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "407.75"

# Row 6 - Solana
$ws.Range("D6").Value = "131.53"
$ws.Range("E6").Value = "  +16.92%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.466.52"
$ws.Range("E7").Value = "  +3.70%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +2.82%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.695"
$ws.Range("E10").Value = "  +8.83%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.132"
$ws.Range("E11").Value = "  +33.22%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "43.55"
$ws.Range("E12").Value = "  +9.67%  "

# Row 13 - now TRON (was WrappedliquidstakedEther2.0)
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.142"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14 - now WrappedliquidstakedEther2.0 (was TRON)
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.023.61"
$ws.Range("E14").Value = "  +3.49%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +5.51%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "20.12"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.480.95"
$ws.Range("E17").Value = "  +4.11%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "63.045.89"
$ws.Range("E18").Value = "  +3.77%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +0.32%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "10.92"
$ws.Range("E20").Value = "  +1.12%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0000143"
$ws.Range("E21").Value = "  +29.05%  "

# Row 22 - ImmutableX
$ws.Range("D22").Value = "3.34"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "82.65"
$ws.Range("E23").Value = "  +10.09%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "13.14"
$ws.Range("E24").Value = "  +1.16%  "

# Row 25 - BitcoinCash
$ws.Range("D25").Value = "312.27"
$ws.Range("E25").Value = "  +3.21%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -0.86%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "30.49"
$ws.Range("E27").Value = "  +6.18%  "

# Row 28 - Filecoin
$ws.Range("E28").Value = "  +2.95%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -1.28%  "

# Row 30 - LEO
$ws.Range("D30").Value = "4.38"
$ws.Range("E30").Value = "  -2.31%  "

# Row 31 - RenderToken
$ws.Range("E31").Value = "  -0.56%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +2.82%  "

# Row 33 - InjectiveProtocol
$ws.Range("D33").Value = "44.29"
$ws.Range("E33").Value = "  +13.83%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "11.84"
$ws.Range("E34").Value = "  +3.77%  "

# Row 35 - Toncoin
$ws.Range("E35").Value = "  -1.21%  "

# Row 36 - Dai
$ws.Range("E36").Value = "  -0.08%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -2.77%  "

# Row 38 - OKB
$ws.Range("D38").Value = "52.63"
$ws.Range("E38").Value = "  +0.38%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +4.80%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.13%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  -2.73%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +2.62%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  +4.20%  "

# Row 44 - Monero
$ws.Range("D44").Value = "137.18"
$ws.Range("E44").Value = "  +0.41%  "

# Row 45 - Celestia
$ws.Range("E45").Value = "  +4.27%  "

# Row 46 - now TheGraph (was NEARProtocol)
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.288"
$ws.Range("E46").Value = "  -3.18%  "

# Row 47 - now NEARProtocol (was TheGraph)
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.98"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  -0.28%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "22.23"
$ws.Range("E49").Value = "  -1.74%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "3.818.32"
$ws.Range("E50").Value = "  +3.44%  "

# Row 51 - Maker
$ws.Range("D51").Value = "2.185.50"
$ws.Range("E51").Value = "  +0.36%  "
